$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$descriptionText = "Entities in this cluster are outgoing and carefree, while being slightly less reserved and imaginative, giving them a relaxed and sociable style. They are outgoing and carefree, making them comfortable and confident in social settings. Their weaknesses are minor being less reserved and imaginative may sometimes lead to less caution, planning, or reflection."
$labelText = "Outgoing and carefree"

$cellA5 = $ws.Range("A5")
$cellB5 = $ws.Range("B5")

$cellA5.Value = $descriptionText
$cellB5.Value = $labelText

$cellA5.WrapText = $true
$cellA5.VerticalAlignment = -4160

$ws.Rows(5).RowHeight = 87

$ws.Range("B5").Select()
